# GPLIM-3541: add Material Type as required header for Manifest uploads
#
# Adds a new "Material Type" column (G) to the manifest worksheet:
#   - G1 header "Material Type", bold white text on a black fill, centered
#   - G2:G24 data cells "DNA:Genomic", centered (reuses the existing
#     center-aligned default style already used elsewhere in the sheet)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell G1 -------------------------------------------------------
$header = $ws.Range("G1")
$header.Value = "Material Type"
$header.HorizontalAlignment = -4108   # xlCenter
$header.Font.Bold = $true
$header.Font.Color = 16777215         # RGB(255,255,255) white
$header.Interior.Color = 0            # RGB(0,0,0) black (foreground)
$header.Interior.PatternColor = 0     # RGB(0,0,0) black (background)

# --- Data cells G2:G24 -----------------------------------------------------
$data = $ws.Range("G2:G24")
$data.Value = "DNA:Genomic"
$data.HorizontalAlignment = -4108     # xlCenter

# --- Selection on the new column -------------------------------------------
$ws.Range("G1:G24").Select()
